# Auto-generated Excel COM-interop script
# Applies scheduled market-data/profit refresh values to the Phantom Profits workbook
# (columns H..N = currentAveragePrice*, LevePrice*, LeveProfit* on sheets ALC/ARM/BSM/CRP/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H55").Value = 393.46155
$ws.Range("I55").Value = 266.83334
$ws.Range("K55").Value = 266.83334
$ws.Range("M55").Value = -52.83334000000002

$ws.Range("H70").Value = 3920.2
$ws.Range("J70").Value = 1500
$ws.Range("L70").Value = 4500
$ws.Range("N70").Value = -5040

$ws.Range("H73").Value = 3920.2
$ws.Range("J73").Value = 1500
$ws.Range("L73").Value = 4500
$ws.Range("N73").Value = -6372

$ws.Range("H74").Value = 4741.5557
$ws.Range("I74").Value = 4084.25
$ws.Range("K74").Value = 4084.25
$ws.Range("M74").Value = -3148.25

$ws.Range("H77").Value = 4741.5557
$ws.Range("I77").Value = 4084.25
$ws.Range("K77").Value = 20421.25
$ws.Range("M77").Value = -15741.25

$ws.Range("H98").Value = 2229.25
$ws.Range("I98").Value = 1207.3334
$ws.Range("J98").Value = 5295
$ws.Range("K98").Value = 1207.3334
$ws.Range("L98").Value = 5295
$ws.Range("M98").Value = 290.6666
$ws.Range("N98").Value = -8291

$ws.Range("H122").Value = 2229.25
$ws.Range("I122").Value = 1207.3334
$ws.Range("J122").Value = 5295
$ws.Range("K122").Value = 3622.0002
$ws.Range("L122").Value = 15885
$ws.Range("M122").Value = -1172.0002
$ws.Range("N122").Value = -20785

$ws.Range("H138").Value = 2376.7778
$ws.Range("J138").Value = 4999.3335
$ws.Range("L138").Value = 14998.0005
$ws.Range("N138").Value = -25278.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12989.947
$ws.Range("I32").Value = 12989.947
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 12989.947
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -12702.947
$ws.Range("N32").ClearContents()

$ws.Range("H61").Value = 1464.6875
$ws.Range("I61").Value = 1481.1428
$ws.Range("J61").Value = 1349.5
$ws.Range("K61").Value = 1481.1428
$ws.Range("L61").Value = 1349.5
$ws.Range("M61").Value = -1269.1428
$ws.Range("N61").Value = -1773.5

$ws.Range("H74").Value = 1797.5
$ws.Range("I74").Value = 1797.5
$ws.Range("K74").Value = 1797.5
$ws.Range("M74").Value = -923.5

$ws.Range("H77").Value = 1797.5
$ws.Range("I77").Value = 1797.5
$ws.Range("K77").Value = 8987.5
$ws.Range("M77").Value = -4619.5

$ws.Range("H102").Value = 1812.9286
$ws.Range("I102").Value = 1782.6666
$ws.Range("K102").Value = 1782.6666
$ws.Range("M102").Value = -160.6666

$ws.Range("H110").Value = 4269.4
$ws.Range("I110").Value = 4530.857
$ws.Range("J110").Value = 609
$ws.Range("K110").Value = 4530.857
$ws.Range("L110").Value = 609
$ws.Range("M110").Value = -2485.857
$ws.Range("N110").Value = -4699

$ws.Range("H122").Value = 2490.6155
$ws.Range("I122").Value = 2490.6155
$ws.Range("K122").Value = 7471.8465
$ws.Range("M122").Value = -5021.8465

$ws.Range("H132").Value = 3521.9443
$ws.Range("I132").Value = 1499.7273
$ws.Range("K132").Value = 4499.1819
$ws.Range("M132").Value = -1969.1819

$ws.Range("H136").Value = 1464.6875
$ws.Range("I136").Value = 1481.1428
$ws.Range("J136").Value = 1349.5
$ws.Range("K136").Value = 4443.428400000001
$ws.Range("L136").Value = 4048.5
$ws.Range("M136").Value = -1893.428400000001
$ws.Range("N136").Value = -9148.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3093.625
$ws.Range("I20").Value = 2351
$ws.Range("K20").Value = 2351
$ws.Range("M20").Value = -2104

$ws.Range("H26").Value = 14331.333
$ws.Range("I26").Value = 14331.333
$ws.Range("K26").Value = 14331.333
$ws.Range("M26").Value = -14039.333

$ws.Range("H75").Value = 4400
$ws.Range("I75").Value = 4400
$ws.Range("K75").Value = 4400
$ws.Range("M75").Value = -3464

$ws.Range("H78").Value = 4400
$ws.Range("I78").Value = 4400
$ws.Range("K78").Value = 13200
$ws.Range("M78").Value = -8520

$ws.Range("H86").Value = 2169.3103
$ws.Range("I86").Value = 2358.7727
$ws.Range("J86").Value = 1573.8572
$ws.Range("K86").Value = 2358.7727
$ws.Range("L86").Value = 1573.8572
$ws.Range("M86").Value = -1235.7727
$ws.Range("N86").Value = -3819.8572

$ws.Range("H89").Value = 2169.3103
$ws.Range("I89").Value = 2358.7727
$ws.Range("J89").Value = 1573.8572
$ws.Range("K89").Value = 11793.8635
$ws.Range("L89").Value = 7869.286
$ws.Range("M89").Value = -6177.863499999999
$ws.Range("N89").Value = -19101.286

$ws.Range("H107").Value = 555.9
$ws.Range("I107").Value = 555.9
$ws.Range("K107").Value = 555.9
$ws.Range("M107").Value = 1364.1

$ws.Range("H134").Value = 3227.2727
$ws.Range("I134").Value = 3227.2727
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9681.8181
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7146.8181
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 37833.668
$ws.Range("J22").Value = 37833.668
$ws.Range("L22").Value = 37833.668
$ws.Range("N22").Value = -38533.668

$ws.Range("H62").Value = 2970.75
$ws.Range("I62").Value = 2947
$ws.Range("K62").Value = 2947
$ws.Range("M62").Value = -2323

$ws.Range("H65").Value = 2970.75
$ws.Range("I65").Value = 2947
$ws.Range("K65").Value = 14735
$ws.Range("M65").Value = -11615

$ws.Range("H122").Value = 1719.7778
$ws.Range("I122").Value = 1747.25
$ws.Range("K122").Value = 5241.75
$ws.Range("M122").Value = -2791.75

$ws.Range("H132").Value = 2505.5
$ws.Range("I132").Value = 2505.5
$ws.Range("K132").Value = 7516.5
$ws.Range("M132").Value = -4986.5

$ws.Range("H134").Value = 3428.6365
$ws.Range("I134").Value = 3428.6365
$ws.Range("K134").Value = 10285.9095
$ws.Range("M134").Value = -7750.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2003.6666
$ws.Range("I102").Value = 1913.091
$ws.Range("K102").Value = 1913.091
$ws.Range("M102").Value = -291.0909999999999

$ws.Range("H122").Value = 3133.889
$ws.Range("I122").Value = 2534.3333
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 7602.999899999999
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -5152.999899999999
$ws.Range("N122").Value = -17899

$ws.Range("H132").Value = 1797.1111
$ws.Range("I132").Value = 1827.1428
$ws.Range("K132").Value = 5481.428400000001
$ws.Range("M132").Value = -2951.428400000001
$ws.Range("N132").Value = -10136

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3337.1667
$ws.Range("I7").Value = 3478.818
$ws.Range("K7").Value = 3478.818
$ws.Range("M7").Value = -3366.818

$ws.Range("H68").Value = 5083.1665
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 5083.1665
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488

$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws.Range("H122").Value = 4521.2856
$ws.Range("I122").Value = 4441.5
$ws.Range("K122").Value = 13324.5
$ws.Range("M122").Value = -10874.5

$ws.Range("H126").Value = 3337.1667
$ws.Range("I126").Value = 3478.818
$ws.Range("K126").Value = 10436.454
$ws.Range("M126").Value = -7966.454000000002

$ws.Range("H132").Value = 3567.923
$ws.Range("I132").Value = 3529.3
$ws.Range("J132").Value = 3696.6667
$ws.Range("K132").Value = 10587.9
$ws.Range("L132").Value = 11090.0001
$ws.Range("M132").Value = -8057.900000000001
$ws.Range("N132").Value = -16150.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 38497.5
$ws.Range("I45").Value = 26247.25
$ws.Range("J45").Value = 50747.75
$ws.Range("K45").Value = 26247.25
$ws.Range("L45").Value = 50747.75
$ws.Range("M45").Value = -25756.25
$ws.Range("N45").Value = -51729.75

$ws.Range("H107").Value = 5008.1113
$ws.Range("I107").Value = 3515.5
$ws.Range("J107").Value = 7993.3335
$ws.Range("K107").Value = 10546.5
$ws.Range("L107").Value = 23980.0005
$ws.Range("M107").Value = -8626.5
$ws.Range("N107").Value = -27820.0005

$ws.Range("H122").Value = 2888.0557
$ws.Range("I122").Value = 2799.0625
$ws.Range("K122").Value = 8397.1875
$ws.Range("M122").Value = -5947.1875

$ws.Range("H132").Value = 2500.762
$ws.Range("I132").Value = 2585.85
$ws.Range("K132").Value = 7757.549999999999
$ws.Range("M132").Value = -5227.549999999999

$ws.Range("H136").Value = 3441.3696
$ws.Range("I136").Value = 3575.8845
$ws.Range("J136").Value = 3266.5
$ws.Range("K136").Value = 10727.6535
$ws.Range("L136").Value = 9799.5
$ws.Range("M136").Value = -8177.6535
$ws.Range("N136").Value = -14899.5
